# Daily BP terminal gate pricing (TGP) rollover.
#
# The source feed advances one business day: each location's previous
# "today" row becomes the new "yesterday" row (values carried forward
# unchanged), and a fresh "today" row is populated with the latest
# effective-date price quotes. Only the date + price cells change; the
# terminal names, labels and styling are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title row grew slightly taller in this save.
$ws.Rows.Item(1).RowHeight = 23.25

# Row 8: Sydney-Botany
$ws.Cells.Item(8,1).Value = 46057
$ws.Cells.Item(8,4).Value = 158.64
$ws.Cells.Item(8,5).Value = 149.27
$ws.Cells.Item(8,6).Value = 159.27
$ws.Cells.Item(8,7).Value = 149.16

# Row 9: Sydney-Silverwater
$ws.Cells.Item(9,1).Value = 46057
$ws.Cells.Item(9,4).Value = 158.64
$ws.Cells.Item(9,5).Value = 149.27
$ws.Cells.Item(9,6).Value = 159.27
$ws.Cells.Item(9,7).Value = 149.16

# Row 10: Newcastle
$ws.Cells.Item(10,1).Value = 46057
$ws.Cells.Item(10,4).Value = 160.33
$ws.Cells.Item(10,5).Value = 150.74
$ws.Cells.Item(10,6).Value = 160.74
$ws.Cells.Item(10,7).Value = 150.99

# Row 11: Sydney-Botany
$ws.Cells.Item(11,1).Value = 46056
$ws.Cells.Item(11,4).Value = 158.47
$ws.Cells.Item(11,5).Value = 149.37
$ws.Cells.Item(11,6).Value = 159.37
$ws.Cells.Item(11,7).Value = 149.26

# Row 12: Sydney-Silverwater
$ws.Cells.Item(12,1).Value = 46056
$ws.Cells.Item(12,4).Value = 158.47
$ws.Cells.Item(12,5).Value = 149.37
$ws.Cells.Item(12,6).Value = 159.37
$ws.Cells.Item(12,7).Value = 149.26

# Row 13: Newcastle
$ws.Cells.Item(13,1).Value = 46056
$ws.Cells.Item(13,4).Value = 159.89
$ws.Cells.Item(13,5).Value = 150.89
$ws.Cells.Item(13,6).Value = 160.89
$ws.Cells.Item(13,7).Value = 151.14

# Row 17: Darwin
$ws.Cells.Item(17,1).Value = 46057
$ws.Cells.Item(17,4).Value = 163.94
$ws.Cells.Item(17,5).Value = 154.2
$ws.Cells.Item(17,6).Value = 164.2

# Row 18: Darwin
$ws.Cells.Item(18,1).Value = 46056
$ws.Cells.Item(18,4).Value = 163.48
$ws.Cells.Item(18,5).Value = 154.35
$ws.Cells.Item(18,6).Value = 164.35

# Row 22: Brisbane
$ws.Cells.Item(22,1).Value = 46057
$ws.Cells.Item(22,4).Value = 159.82
$ws.Cells.Item(22,5).Value = 150.73
$ws.Cells.Item(22,6).Value = 160.33
$ws.Cells.Item(22,7).Value = 152.49

# Row 23: Cairns
$ws.Cells.Item(23,1).Value = 46057
$ws.Cells.Item(23,4).Value = 165.11
$ws.Cells.Item(23,5).Value = 156.86
$ws.Cells.Item(23,6).Value = 166.86

# Row 24: Gladstone
$ws.Cells.Item(24,1).Value = 46057
$ws.Cells.Item(24,4).Value = 165.27
$ws.Cells.Item(24,5).Value = 157.5
$ws.Cells.Item(24,6).Value = 167.5

# Row 25: Mackay
$ws.Cells.Item(25,1).Value = 46057
$ws.Cells.Item(25,4).Value = 165.27
$ws.Cells.Item(25,5).Value = 157.03
$ws.Cells.Item(25,6).Value = 167.03
$ws.Cells.Item(25,7).Value = 157.88

# Row 26: Townsville
$ws.Cells.Item(26,1).Value = 46057
$ws.Cells.Item(26,4).Value = 164.86
$ws.Cells.Item(26,5).Value = 158.61
$ws.Cells.Item(26,6).Value = 168.61

# Row 27: Brisbane
$ws.Cells.Item(27,1).Value = 46056
$ws.Cells.Item(27,4).Value = 159.54
$ws.Cells.Item(27,5).Value = 150.83
$ws.Cells.Item(27,6).Value = 160.43
$ws.Cells.Item(27,7).Value = 152.58

# Row 28: Cairns
$ws.Cells.Item(28,1).Value = 46056
$ws.Cells.Item(28,4).Value = 164.66
$ws.Cells.Item(28,5).Value = 157.01
$ws.Cells.Item(28,6).Value = 167.01

# Row 29: Gladstone
$ws.Cells.Item(29,1).Value = 46056
$ws.Cells.Item(29,4).Value = 164.82
$ws.Cells.Item(29,5).Value = 157.65
$ws.Cells.Item(29,6).Value = 167.65

# Row 30: Mackay
$ws.Cells.Item(30,1).Value = 46056
$ws.Cells.Item(30,4).Value = 164.82
$ws.Cells.Item(30,5).Value = 157.18
$ws.Cells.Item(30,6).Value = 167.18
$ws.Cells.Item(30,7).Value = 158.03

# Row 31: Townsville
$ws.Cells.Item(31,1).Value = 46056
$ws.Cells.Item(31,4).Value = 164.4
$ws.Cells.Item(31,5).Value = 158.76
$ws.Cells.Item(31,6).Value = 168.76

# Row 35: Adelaide
$ws.Cells.Item(35,1).Value = 46057
$ws.Cells.Item(35,4).Value = 158.83
$ws.Cells.Item(35,5).Value = 148.54
$ws.Cells.Item(35,6).Value = 157.54

# Row 36: Adelaide
$ws.Cells.Item(36,1).Value = 46056
$ws.Cells.Item(36,4).Value = 158.82
$ws.Cells.Item(36,5).Value = 148.69
$ws.Cells.Item(36,6).Value = 157.69

# Row 40: Burnie
$ws.Cells.Item(40,1).Value = 46057
$ws.Cells.Item(40,4).Value = 164.65
$ws.Cells.Item(40,5).Value = 156.04
$ws.Cells.Item(40,6).Value = 166.04

# Row 41: Hobart
$ws.Cells.Item(41,1).Value = 46057
$ws.Cells.Item(41,4).Value = 164.36
$ws.Cells.Item(41,5).Value = 156.46
$ws.Cells.Item(41,6).Value = 166.46

# Row 42: Burnie
$ws.Cells.Item(42,1).Value = 46056
$ws.Cells.Item(42,4).Value = 164.2
$ws.Cells.Item(42,5).Value = 156.07
$ws.Cells.Item(42,6).Value = 166.07

# Row 43: Hobart
$ws.Cells.Item(43,1).Value = 46056
$ws.Cells.Item(43,4).Value = 163.92
$ws.Cells.Item(43,5).Value = 156.49
$ws.Cells.Item(43,6).Value = 166.49

# Row 47: Geelong
$ws.Cells.Item(47,1).Value = 46057
$ws.Cells.Item(47,4).Value = 158.93
$ws.Cells.Item(47,5).Value = 150.22
$ws.Cells.Item(47,6).Value = 160.22

# Row 48: Melbourne
$ws.Cells.Item(48,1).Value = 46057
$ws.Cells.Item(48,4).Value = 158.55
$ws.Cells.Item(48,5).Value = 150.16
$ws.Cells.Item(48,6).Value = 160.16

# Row 49: Geelong
$ws.Cells.Item(49,1).Value = 46056
$ws.Cells.Item(49,4).Value = 157.99
$ws.Cells.Item(49,5).Value = 150.68
$ws.Cells.Item(49,6).Value = 160.68

# Row 50: Melbourne
$ws.Cells.Item(50,1).Value = 46056
$ws.Cells.Item(50,4).Value = 157.61
$ws.Cells.Item(50,5).Value = 150.62
$ws.Cells.Item(50,6).Value = 160.62

# Row 54: Broome
$ws.Cells.Item(54,1).Value = 46057
$ws.Cells.Item(54,4).Value = 173.65
$ws.Cells.Item(54,5).Value = 164.28
$ws.Cells.Item(54,6).Value = 174.28

# Row 55: Esperance
$ws.Cells.Item(55,1).Value = 46057
$ws.Cells.Item(55,4).Value = 163.07
$ws.Cells.Item(55,5).Value = 162.12
$ws.Cells.Item(55,6).Value = 172.12

# Row 56: Geraldton
$ws.Cells.Item(56,1).Value = 46057
$ws.Cells.Item(56,4).Value = 162.95

# Row 57: Kalgoorlie
$ws.Cells.Item(57,1).Value = 46057
$ws.Cells.Item(57,4).Value = 163.53
$ws.Cells.Item(57,5).Value = 156.54

# Row 58: Perth
$ws.Cells.Item(58,1).Value = 46057
$ws.Cells.Item(58,4).Value = 159.3
$ws.Cells.Item(58,5).Value = 152.44
$ws.Cells.Item(58,6).Value = 162.44

# Row 59: Port Hedland
$ws.Cells.Item(59,1).Value = 46057
$ws.Cells.Item(59,4).Value = 166.32
$ws.Cells.Item(59,5).Value = 162.48

# Row 60: Broome
$ws.Cells.Item(60,1).Value = 46056
$ws.Cells.Item(60,4).Value = 173.19
$ws.Cells.Item(60,5).Value = 164.46
$ws.Cells.Item(60,6).Value = 174.46

# Row 61: Esperance
$ws.Cells.Item(61,1).Value = 46056
$ws.Cells.Item(61,4).Value = 162.62
$ws.Cells.Item(61,5).Value = 162.26
$ws.Cells.Item(61,6).Value = 172.26

# Row 62: Geraldton
$ws.Cells.Item(62,1).Value = 46056
$ws.Cells.Item(62,4).Value = 162.5

# Row 63: Kalgoorlie
$ws.Cells.Item(63,1).Value = 46056
$ws.Cells.Item(63,4).Value = 163.08
$ws.Cells.Item(63,5).Value = 156.68

# Row 64: Perth
$ws.Cells.Item(64,1).Value = 46056
$ws.Cells.Item(64,4).Value = 158.85
$ws.Cells.Item(64,5).Value = 152.58
$ws.Cells.Item(64,6).Value = 162.58

# Row 65: Port Hedland
$ws.Cells.Item(65,1).Value = 46056
$ws.Cells.Item(65,4).Value = 165.86
$ws.Cells.Item(65,5).Value = 162.64
